$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the mapping value for "Periodicity" field key: field_wbddh_periodicity -> field_frequency
$ws.Range("B10").Value = "field_frequency"

# Match the cursor/selection position recorded in the saved file (E9)
$ws.Range("E9").Select()
